$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => (D, M, N, O, P, S) target values after the weekly shuffle/update
$data = @(
    ,@(2, 44434, 100, 20000, 21000, 20500, 1025)
    ,@(3, 44333, 100, 19500, 20000, 19750, 988)
    ,@(4, 44466, 100, 20000, 21000, 20500, 1025)
    ,@(5, 44301, 100, 18000, 19000, 18500, 925)
    ,@(6, 44343, 100, 19500, 20000, 19750, 988)
    ,@(7, 44467, 200, 20000, 21000, 20500, 1025)
    ,@(8, 44448, 100, 20000, 21000, 20500, 1025)
    ,@(9, 44364, 140, 20000, 21000, 20500, 1025)
    ,@(10, 44442, 140, 20000, 21000, 20500, 1025)
    ,@(11, 44410, 200, 20000, 21000, 20500, 1025)
    ,@(12, 44365, 100, 20000, 21000, 20500, 1025)
    ,@(13, 44420, 160, 20000, 21000, 20500, 1025)
    ,@(14, 44431, 160, 21000, 22000, 21500, 1075)
    ,@(15, 44326, 160, 19500, 20000, 19750, 988)
    ,@(16, 44336, 100, 19500, 20000, 19750, 988)
    ,@(17, 44441, 160, 20000, 21000, 20500, 1025)
    ,@(18, 44407, 160, 20000, 21000, 20500, 1025)
    ,@(19, 44417, 160, 20000, 21000, 20500, 1025)
    ,@(20, 44428, 100, 20000, 21000, 20500, 1025)
    ,@(21, 44315, 100, 20000, 21000, 20500, 1025)
    ,@(22, 44462, 100, 19500, 20000, 19750, 988)
    ,@(23, 44435, 260, 20000, 22000, 21115, 1056)
    ,@(24, 44473, 40, 19500, 20000, 19750, 988)
    ,@(25, 44445, 160, 20000, 21000, 20500, 1025)
    ,@(26, 44427, 200, 20000, 21000, 20500, 1025)
    ,@(27, 44418, 200, 20000, 21000, 20500, 1025)
    ,@(28, 44350, 160, 19000, 20000, 19500, 975)
    ,@(29, 44335, 200, 19000, 20000, 19500, 975)
    ,@(30, 44474, 200, 19000, 20000, 19500, 975)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value  = $row[1]   # D: Fecha
    $ws.Cells.Item($r, 13).Value = $row[2]   # M: Volumen
    $ws.Cells.Item($r, 14).Value = $row[3]   # N: Precio minimo
    $ws.Cells.Item($r, 15).Value = $row[4]   # O: Precio maximo
    $ws.Cells.Item($r, 16).Value = $row[5]   # P: Precio promedio ponderado
    $ws.Cells.Item($r, 19).Value = $row[6]   # S: Precio $/Kg
}
